# Remove the first six paragraphs (the "039/040 ... youtube links" block)
# from the "Text Placeholder 2" shape on slide 4, leaving only the
# trailing empty paragraph that was already present after them.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(4)
$tr = $shp.TextFrame.TextRange

# Repeatedly delete the first paragraph six times; each deletion shifts
# the following paragraphs up, so paragraph 1 is always the next one to
# remove. This leaves the original trailing (7th) paragraph untouched.
for ($i = 1; $i -le 6; $i++) {
    $para = $tr.Paragraphs(1, 1)
    $para.Delete()
}
